$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 3) for the S&P500 index, mirroring the existing
# Security Label / Ticker / ISIN / Is Index / Last Price columns (A:E).
$ws.Range("A3").Value = "S&P500"
$ws.Range("B3").Value = "SPX"
$ws.Range("C3").Value = 0

# "Is Index" needs to hold the literal text "True" (matching the existing
# "False" text cell above it) rather than Excel's native boolean TRUE.
# Writing the word directly makes Excel auto-coerce it to a boolean, so
# build it as a text formula first and then flatten it back down to a
# plain value via copy / paste-special, which keeps it a genuine
# text/shared-string cell instead of a boolean or formula cell.
$ws.Range("D3").Formula = "=""True"""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)

# Last Price for the S&P500 row.
$ws.Range("E3").Value = 5768
